$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: insert a new "sim" cell at H20, copying H2's format (style) so the
# new cell matches the existing "sim" cells (style index 1) instead of
# picking up the sheet's default (unstyled) format.
$ws.Range("H20").Value = "sim"
$ws.Range("H2").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 7: the "sim" flag moves from J7 to I7 (I7 was previously blank).
$ws.Range("I7").Value = "sim"
$ws.Range("J7").Clear()

# Row 2: drop the extra "sim" cell at J2 entirely.
$ws.Range("J2").Clear()
